# Update LR-pair sheet with new TPM-derived expression values and
# recompute all dependent (specificity / edge-weight) columns, mirroring
# the NATMI LR-pairs pipeline recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand average expression values (TPM) per sending cluster (column G)
$ligandAvg = @{
    "ECs"            = 45.89896166666667
    "FAPs"           = 3.815058666666667
    "MuSCs"          = 4.651706333333334
    "Resolving-Mac"  = 26.182385
}

# New receptor average expression values (TPM) per target cluster (column M)
$receptorAvg = @{
    "ECs"            = 168.1098273333333
    "FAPs"           = 163.0062356666667
    "MuSCs"          = 165.99353
    "Resolving-Mac"  = 66.22673433333334
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp

# First pass: write the new average-expression base values (G, M) and
# their totals (H, N) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2   # column A
    $target  = $ws.Cells.Item($r, 4).Value2   # column D

    $nExprLigand   = $ws.Cells.Item($r, 5).Value2   # column E
    $nExprReceptor = $ws.Cells.Item($r, 11).Value2  # column K

    $g = $ligandAvg[$sending]
    $m = $receptorAvg[$target]

    $ws.Cells.Item($r, 7).Value2  = $g                  # G - ligand average expression value
    $ws.Cells.Item($r, 8).Value2  = $g * $nExprLigand   # H - ligand total expression value

    $ws.Cells.Item($r, 13).Value2 = $m                  # M - receptor average expression value
    $ws.Cells.Item($r, 14).Value2 = $m * $nExprReceptor # N - receptor total expression value
}

# Sums used to normalise the specificity columns
$sumLigandAvg = 0.0
foreach ($v in $ligandAvg.Values) { $sumLigandAvg += $v }

$sumReceptorAvg = 0.0
foreach ($v in $receptorAvg.Values) { $sumReceptorAvg += $v }

# Second pass: recompute the derived specificity / edge-weight columns.
for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2    # G
    $h = $ws.Cells.Item($r, 8).Value2    # H
    $m = $ws.Cells.Item($r, 13).Value2   # M
    $n = $ws.Cells.Item($r, 14).Value2   # N

    $specLigand   = $g / $sumLigandAvg
    $specReceptor = $m / $sumReceptorAvg

    $ws.Cells.Item($r, 9).Value2  = $specLigand     # I - ligand specificity (avg)
    $ws.Cells.Item($r, 10).Value2 = $specLigand     # J - ligand specificity (total)

    $ws.Cells.Item($r, 15).Value2 = $specReceptor   # O - receptor specificity (avg)
    $ws.Cells.Item($r, 16).Value2 = $specReceptor   # P - receptor specificity (total)

    $edgeAvg   = $g * $m
    $edgeTotal = $h * $n
    $edgeSpec  = $specLigand * $specReceptor

    $ws.Cells.Item($r, 17).Value2 = $edgeAvg        # Q - edge average expression weight
    $ws.Cells.Item($r, 18).Value2 = $edgeTotal      # R - edge total expression weight
    $ws.Cells.Item($r, 19).Value2 = $edgeSpec       # S - edge average expression specificity
    $ws.Cells.Item($r, 20).Value2 = $edgeSpec       # T - edge total expression specificity
}
